$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update simulation-derived probability matrix values (team_specific_matrix)
# as re-computed after adding more simulated games.
    $ws.Range("B2").Value = 0.2269230769230769
    $ws.Range("C2").Value = 0.5
    $ws.Range("J2").Value = 0.01153846153846154
    $ws.Range("P2").Value = 0.15
    $ws.Range("S2").Value = 0.1115384615384615
    $ws.Range("B3").Value = 0.02205882352941177
    $ws.Range("C3").Value = 0.03676470588235294
    $ws.Range("J3").Value = 0.01470588235294118
    $ws.Range("P3").Value = 0.75
    $ws.Range("S3").Value = 0.1764705882352941
    $ws.Range("J4").Value = 0.09090909090909091
    $ws.Range("P4").Value = 0.7878787878787878
    $ws.Range("S4").Value = 0.1212121212121212
    $ws.Range("B6").Value = 0.09895833333333333
    $ws.Range("D6").Value = 0.005208333333333333
    $ws.Range("F6").Value = 0.0625
    $ws.Range("J6").Value = 0.3177083333333333
    $ws.Range("Q6").Value = 0.1145833333333333
    $ws.Range("R6").Value = 0.04166666666666666
    $ws.Range("S6").Value = 0.359375
    $ws.Range("B7").Value = 0.08260869565217391
    $ws.Range("D7").Value = 0.01739130434782609
    $ws.Range("F7").Value = 0.03478260869565217
    $ws.Range("J7").Value = 0.1695652173913043
    $ws.Range("O7").Value = 0.01739130434782609
    $ws.Range("Q7").Value = 0.1260869565217391
    $ws.Range("R7").Value = 0.1173913043478261
    $ws.Range("S7").Value = 0.4347826086956522
    $ws.Range("B8").Value = 0.07597535934291581
    $ws.Range("D8").Value = 0.01026694045174538
    $ws.Range("E8").Value = 0.002053388090349076
    $ws.Range("F8").Value = 0.04106776180698152
    $ws.Range("J8").Value = 0.09856262833675565
    $ws.Range("O8").Value = 0.02464065708418891
    $ws.Range("Q8").Value = 0.1211498973305955
    $ws.Range("R8").Value = 0.1416837782340862
    $ws.Range("S8").Value = 0.484599589322382
    $ws.Range("B9").Value = 0.0915032679738562
    $ws.Range("D9").Value = 0.0130718954248366
    $ws.Range("F9").Value = 0.05228758169934641
    $ws.Range("J9").Value = 0.0915032679738562
    $ws.Range("O9").Value = 0.006535947712418301
    $ws.Range("Q9").Value = 0.1372549019607843
    $ws.Range("R9").Value = 0.1241830065359477
    $ws.Range("S9").Value = 0.4836601307189543
    $ws.Range("B10").Value = 0.08498349834983498
    $ws.Range("D10").Value = 0.0165016501650165
    $ws.Range("E10").Value = 0.0008250825082508251
    $ws.Range("F10").Value = 0.05693069306930693
    $ws.Range("J10").Value = 0.1047854785478548
    $ws.Range("O10").Value = 0.01155115511551155
    $ws.Range("Q10").Value = 0.2161716171617162
    $ws.Range("R10").Value = 0.1014851485148515
    $ws.Range("S10").Value = 0.4067656765676568
    $ws.Range("G11").Value = 0.1448863636363636
    $ws.Range("J11").Value = 0.08522727272727272
    $ws.Range("K11").Value = 0.2130681818181818
    $ws.Range("L11").Value = 0.5340909090909091
    $ws.Range("S11").Value = 0.02272727272727273
    $ws.Range("F12").Value = 0.005181347150259068
    $ws.Range("G12").Value = 0.7875647668393783
    $ws.Range("J12").Value = 0.155440414507772
    $ws.Range("K12").Value = 0.01036269430051814
    $ws.Range("L12").Value = 0.0310880829015544
    $ws.Range("S12").Value = 0.01036269430051814
    $ws.Range("G13").Value = 0.8
    $ws.Range("J13").Value = 0.15
    $ws.Range("S13").Value = 0.05
    $ws.Range("F15").Value = 0.03317535545023697
    $ws.Range("H15").Value = 0.1848341232227488
    $ws.Range("I15").Value = 0.07109004739336493
    $ws.Range("J15").Value = 0.2938388625592417
    $ws.Range("K15").Value = 0.07582938388625593
    $ws.Range("M15").Value = 0.01421800947867299
    $ws.Range("N15").Value = 0.004739336492890996
    $ws.Range("O15").Value = 0.1090047393364929
    $ws.Range("S15").Value = 0.2132701421800948
    $ws.Range("F16").Value = 0.03680981595092025
    $ws.Range("H16").Value = 0.1717791411042945
    $ws.Range("I16").Value = 0.1104294478527607
    $ws.Range("J16").Value = 0.4294478527607362
    $ws.Range("K16").Value = 0.0736196319018405
    $ws.Range("M16").Value = 0.01226993865030675
    $ws.Range("O16").Value = 0.04294478527607362
    $ws.Range("S16").Value = 0.1226993865030675
    $ws.Range("F17").Value = 0.01794871794871795
    $ws.Range("H17").Value = 0.1871794871794872
    $ws.Range("I17").Value = 0.05641025641025641
    $ws.Range("J17").Value = 0.4282051282051282
    $ws.Range("K17").Value = 0.1025641025641026
    $ws.Range("M17").Value = 0.02307692307692308
    $ws.Range("O17").Value = 0.07179487179487179
    $ws.Range("S17").Value = 0.1128205128205128
    $ws.Range("F18").Value = 0.02845528455284553
    $ws.Range("H18").Value = 0.1585365853658537
    $ws.Range("I18").Value = 0.06504065040650407
    $ws.Range("J18").Value = 0.4471544715447154
    $ws.Range("K18").Value = 0.1260162601626016
    $ws.Range("M18").Value = 0.004065040650406504
    $ws.Range("O18").Value = 0.07317073170731707
    $ws.Range("S18").Value = 0.0975609756097561
    $ws.Range("F19").Value = 0.01657874905802562
    $ws.Range("H19").Value = 0.2351168048229088
    $ws.Range("I19").Value = 0.06405425772418991
    $ws.Range("J19").Value = 0.346646571213263
    $ws.Range("K19").Value = 0.1311228334589299
    $ws.Range("M19").Value = 0.01959306706857573
    $ws.Range("N19").Value = 0.001507159005275057
    $ws.Range("O19").Value = 0.05802562170308968
    $ws.Range("S19").Value = 0.1273549359457423
